# Append: 2025-09-23 06:33 JST
# Update the "取得日時" (acquired timestamp) column (A) for the existing
# rows on the "ランサーズ" sheet from 2025-09-23 06:27:21 to 2025-09-23 06:33:44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldTimestamp = "2025-09-23 06:27:21"
$newTimestamp = "2025-09-23 06:33:44"

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
